# Loan RBI, Variable Instalments
# - Insert a new (blank) column before column N on the "Repayment schedule"
#   sheet, shifting the old N/O/P ("Late" / blank-heading / "Outstanding")
#   columns one place to the right.
# - Update the active selection on "Edit Repayment Schedule" and
#   "Repayment schedule" sheets.

$wb = $excel.ActiveWorkbook

$wsEditSchedule = $wb.Worksheets.Item("Edit Repayment Schedule")
$wsSchedule     = $wb.Worksheets.Item("Repayment schedule")

# Insert a new blank column at N (pushes Late/blank/Outstanding to O/P/Q).
$wsSchedule.Columns("N").Insert()

# The newly inserted column keeps the same width as its left neighbour (M),
# i.e. an OOXML column width of 11 characters.
$wsSchedule.Columns("N").ColumnWidth = 10.17

# Update the selected cell on each sheet to match the saved view state.
$wsSchedule.Activate() | Out-Null
$wsSchedule.Range("J19").Select() | Out-Null

$wsEditSchedule.Activate() | Out-Null
$wsEditSchedule.Range("E9").Select() | Out-Null
